$d = $word.ActiveDocument

# Paragraph 3 ("Dokumentacja do X10 włącznie ze schematem podłączeń") gets the
# question that used to be in the next paragraph.
$p3 = $d.Paragraphs.Item(3).Range
$p3.Text = "Dlaczego program się wypiedala ? To nie ma prawa się dziać ?. Program była aktualizowany tak aby zasłony mogły być zasłonięte w okresie letnim. A w okresie zimowym mają być cały czas otwarte ?"

# Paragraph 4 (previously holding the text above) becomes a new question.
$p4 = $d.Paragraphs.Item(4).Range
$p4.Text = "Dlaczego nie mogę zrobić uploadu programu z termostatu ?"

# Insert a brand-new list item right after paragraph 4 with another new question.
$p4after = $d.Paragraphs.Item(4).Range
$p4after.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(5).Range
$newPara.Text = "Gdzie znajdę dokumentację do Dialog Box ?"
